$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 565, pushing existing rows 565-681 down to 566-682
$ws.Rows.Item(565).Insert()

# Fill in the new record's data (same fixed fields as the rest of the table)
$ws.Range("A565").Value = 3
$ws.Range("B565").Value = "Femacal de La Calera"
$ws.Range("C565").Value = "Coquimbo"
$ws.Range("D565").Value = 45275
$ws.Range("E565").Value = 5
$ws.Range("F565").Value = 100112040
$ws.Range("G565").Value = "Cilantro"
$ws.Range("H565").Value = "Sin especificar"
$ws.Range("I565").Value = "Primera"
$ws.Range("J565").Value = 220
$ws.Range("K565").Value = 5500
$ws.Range("L565").Value = 6000
$ws.Range("M565").Value = 5750
$ws.Range("N565").Value = "$/docena de atados (3 kilos)"
$ws.Range("O565").Value = "Provincia de Quillota"
$ws.Range("P565").Value = 1917
$ws.Range("Q565").Value = 3
$ws.Range("R565").Value = "Hortaliza"
